# Root_dmg_20.xlsx — standardize variable names / translate headers & sheet title
#
# Summary of the underlying commit:
#  * Sheet renamed from Spanish "Daño en raiz (20 repl.)" to English
#    "Root dmg (20 reps.)".
#  * The table ("Tabla535") header row / column names are translated:
#       Fecha              -> Date
#       ID Parcela         -> Field
#       Tratamiento        -> Treatment
#       Repetición         -> Repeat
#       Peso biomasa raiz  -> Root_weight
#       Observaciones      -> Observations
#  * Column widths for B:F are narrowed/adjusted (no longer "best fit").
#  * Active selection moved to H17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet / tab.
$ws.Name = "Root dmg (20 reps.)"

# 2) Rename the table header cells (this also renames the ListObject's
#    ListColumns, since the table header row drives the column names).
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("D1").Value = "Repeat"
$ws.Range("E1").Value = "Root_weight"
$ws.Range("F1").Value = "Observations"

# 3) Resize columns B:F (no longer auto "best fit" — explicit widths).
$ws.Columns.Item(2).ColumnWidth = 8.6666666666667
$ws.Columns.Item(3).ColumnWidth = 12.1666666666667
$ws.Columns.Item(4).ColumnWidth = 9.3333333333333
$ws.Columns.Item(5).ColumnWidth = 16.6666666666667
$ws.Columns.Item(6).ColumnWidth = 93

# 4) Move the active selection to H17.
$ws.Range("H17").Select() | Out-Null
